$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 12191
$ws1.Range("F10").Value = 2603
$ws1.Range("F11").Value = 1129
$ws1.Range("F12").Value = 200
$ws1.Range("F14").Value = 5293
$ws1.Range("F18").Value = 11457
$ws1.Range("F19").Value = 11532
$ws1.Range("F21").Value = 60

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12191
$ws4.Range("F10").Value = 2603
$ws4.Range("F12").Value = 1129
$ws4.Range("F13").Value = 200
$ws4.Range("F15").Value = 5293
$ws4.Range("F19").Value = 11457
$ws4.Range("F20").Value = 11532
$ws4.Range("F22").Value = 60
